$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "Test - Test"
$ws.Range("D4").Value = "Test - Test"

$ws.Range("E4").Select()
